# The commit swaps the contents of ppt/theme/theme1.xml (the theme attached
# to the slide master / deck design, currently "Integral") and
# ppt/theme/theme2.xml (the theme attached only to the notes master,
# currently "Office Theme") - i.e. theme1.xml becomes the Office Theme
# palette and theme2.xml becomes the Integral palette. In this deck the two
# themes already share an identical font scheme and format scheme (fills,
# lines, effects), so the only substantive difference is the 12-colour
# scheme (and the cosmetic theme/clrScheme "name" attributes, which the
# PowerPoint object model does not expose for editing).
#
# The only theme surface reachable from the PowerPoint object model is the
# one and only design's ThemeColorScheme (SlideMaster.Theme /
# NotesMaster.Theme both resolve to this same theme, which is backed by
# theme1.xml) - so we repoint its 12 colours from the Integral palette to
# the Office Theme palette.

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$colorScheme = $master.Theme.ThemeColorScheme

# Target "Office Theme" colour scheme, in ThemeColorScheme.Item() index
# order: 1 dk1, 2 lt1, 3 dk2, 4 lt2, 5 accent1, 6 accent2, 7 accent3,
# 8 accent4, 9 accent5, 10 accent6, 11 hlink, 12 folHlink.
$officeThemeHex = @(
    "000000",
    "FFFFFF",
    "44546A",
    "E7E6E6",
    "5B9BD5",
    "ED7D31",
    "A5A5A5",
    "FFC000",
    "4472C4",
    "70AD47",
    "0563C1",
    "954F72"
)

for ($i = 1; $i -le $officeThemeHex.Length; $i++) {
    $hex = $officeThemeHex[$i - 1]
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    # PowerPoint RGB colours are stored as OLE_COLOR (0x00BBGGRR).
    $oleColor = ($b * 65536) + ($g * 256) + $r
    $colorScheme.Item($i).RGB = $oleColor
}
